$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# New vintage column BH ("Agosto.2021"), the 60th column.
# ---------------------------------------------------------------------

# Header BH1: same (bold/centered/bordered) formatting as the rest of row 1.
$ws.Range("BG1").Copy()
$ws.Range("BH1").PasteSpecial(-4122)
$ws.Range("BH1").Value = "Agosto.2021"
$excel.CutCopyMode = 0

# Data rows 2-73: the new vintage repeats the previous (BG) vintage's value.
for ($r = 2; $r -le 73; $r++) {
    $v = $ws.Cells.Item($r, 59).Value()
    $ws.Cells.Item($r, 60).Value = $v
}

# Row 74 (period 01-01-2021) gets a revised figure under the new vintage.
$ws.Cells.Item(74, 60).Value = 38076

# ---------------------------------------------------------------------
# New row 75: period 01-04-2021, first reported only under the new vintage.
# ---------------------------------------------------------------------

# Column A holds plain text dates (e.g. "01-01-2021"); entering that string
# directly would get auto-parsed into a date serial by Excel. Stage it in a
# faraway scratch cell formatted as Text, then copy/paste just the value
# across so the destination cell keeps the default (unformatted) style,
# same as the existing date cells in column A.
$scratch = $ws.Cells.Item(1000, 1000)
$scratch.NumberFormat = "@"
$scratch.Value = "01-04-2021"
$scratch.Copy()
$ws.Cells.Item(75, 1).PasteSpecial(-4163)
$scratch.Clear()
$excel.CutCopyMode = 0

$ws.Cells.Item(75, 60).Value = 39677
